# Apply crypto price/symbol updates for 2022-12-18 GitHub Actions run.
# All target cells in this sheet store values as text (inlineStr), including
# numeric-looking prices, so we must write them as text and strip any
# quote-prefix styling Excel may apply, to keep cells unstyled like the source.

function Set-CellText($sheet, $ref, $val) {
    $cell = $sheet.Range($ref)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws "D2" '247.04'
Set-CellText $ws "D4" '5.494'
Set-CellText $ws "D5" '0.05635'
Set-CellText $ws "D7" '0.8066'
Set-CellText $ws "D8" '1.048'
Set-CellText $ws "D9" '0.1445'
Set-CellText $ws "D10" '0.07369'
Set-CellText $ws "D11" '0.03189'
Set-CellText $ws "D12" '0.02934'
Set-CellText $ws "D13" '0.09270'
Set-CellText $ws "D14" '0.001676'
Set-CellText $ws "D15" '3.204'
Set-CellText $ws "D16" '0.04731'
Set-CellText $ws "D17" '0.0005852'
Set-CellText $ws "D18" '0.006288'
Set-CellText $ws "D19" '0.001055'
Set-CellText $ws "D20" '0.004111'
Set-CellText $ws "D22" '3.977'
Set-CellText $ws "D23" '3.386'
Set-CellText $ws "D24" '2.134'
Set-CellText $ws "D25" '0.3274'
Set-CellText $ws "D27" '0.0003012'
Set-CellText $ws "D40" '0.04155'
Set-CellText $ws "B41" 'BKEXToken'
Set-CellText $ws "C41" 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-CellText $ws "D41" '0.1041'
Set-CellText $ws "E41" '40BKEXTokenBKK'
Set-CellText $ws "B42" 'CEJI'
Set-CellText $ws "C42" 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-CellText $ws "D42" '0.003204'
Set-CellText $ws "E42" '41CEJICEJI'
Set-CellText $ws "B43" 'KickToken'
Set-CellText $ws "C43" 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-CellText $ws "D43" '0.003252'
Set-CellText $ws "E43" '42KickTokenKICKWorstin24h'
Set-CellText $ws "D44" '0.009078'
Set-CellText $ws "D45" '0.00005663'
Set-CellText $ws "D46" '0.00000000753'
Set-CellText $ws "D47" '0.6826'
Set-CellText $ws "D48" '0.02040'
Set-CellText $ws "E48" '47BOLOBOLO'
Set-CellText $ws "D49" '0.00002108'
Set-CellText $ws "D50" '0.01014'
